$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 3 new rows before row 224. This pushes the existing rows
# 224-337 down to 227-340 (matching the diff's net "shift by 3" pattern).
$ws.Rows.Item(224).Insert()
$ws.Rows.Item(224).Insert()
$ws.Rows.Item(224).Insert()

# Fill in the 3 newly inserted rows (224-226) with the new weekly price
# record for date 44609 (2022-02-17), mirroring the constant columns
# (A,B,C,E-L,Q,R,T) used throughout this block and the new M/N/O/P/S
# values from the diff.

# Row 224: Especial
$ws.Cells.Item(224, 1).Value = 2
$ws.Cells.Item(224, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(224, 3).Value = "Coquimbo"
$ws.Cells.Item(224, 4).Value = 44609
$ws.Cells.Item(224, 5).Value = 4
$ws.Cells.Item(224, 6).Value = "Fruta"
$ws.Cells.Item(224, 7).Value = 100101
$ws.Cells.Item(224, 8).Value = "Berries"
$ws.Cells.Item(224, 9).Value = 100112025
$ws.Cells.Item(224, 10).Value = "Frutilla"
$ws.Cells.Item(224, 11).Value = "Sin especificar"
$ws.Cells.Item(224, 12).Value = "Especial"
$ws.Cells.Item(224, 13).Value = 400
$ws.Cells.Item(224, 14).Value = 11500
$ws.Cells.Item(224, 15).Value = 12000
$ws.Cells.Item(224, 16).Value = 11750
$ws.Cells.Item(224, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(224, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(224, 19).Value = 1679
$ws.Cells.Item(224, 20).Value = 7

# Row 225: Primera
$ws.Cells.Item(225, 1).Value = 2
$ws.Cells.Item(225, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(225, 3).Value = "Coquimbo"
$ws.Cells.Item(225, 4).Value = 44609
$ws.Cells.Item(225, 5).Value = 4
$ws.Cells.Item(225, 6).Value = "Fruta"
$ws.Cells.Item(225, 7).Value = 100101
$ws.Cells.Item(225, 8).Value = "Berries"
$ws.Cells.Item(225, 9).Value = 100112025
$ws.Cells.Item(225, 10).Value = "Frutilla"
$ws.Cells.Item(225, 11).Value = "Sin especificar"
$ws.Cells.Item(225, 12).Value = "Primera"
$ws.Cells.Item(225, 13).Value = 500
$ws.Cells.Item(225, 14).Value = 9500
$ws.Cells.Item(225, 15).Value = 10000
$ws.Cells.Item(225, 16).Value = 9750
$ws.Cells.Item(225, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(225, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(225, 19).Value = 1393
$ws.Cells.Item(225, 20).Value = 7

# Row 226: Segunda
$ws.Cells.Item(226, 1).Value = 2
$ws.Cells.Item(226, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(226, 3).Value = "Coquimbo"
$ws.Cells.Item(226, 4).Value = 44609
$ws.Cells.Item(226, 5).Value = 4
$ws.Cells.Item(226, 6).Value = "Fruta"
$ws.Cells.Item(226, 7).Value = 100101
$ws.Cells.Item(226, 8).Value = "Berries"
$ws.Cells.Item(226, 9).Value = 100112025
$ws.Cells.Item(226, 10).Value = "Frutilla"
$ws.Cells.Item(226, 11).Value = "Sin especificar"
$ws.Cells.Item(226, 12).Value = "Segunda"
$ws.Cells.Item(226, 13).Value = 400
$ws.Cells.Item(226, 14).Value = 7500
$ws.Cells.Item(226, 15).Value = 8000
$ws.Cells.Item(226, 16).Value = 7750
$ws.Cells.Item(226, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(226, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(226, 19).Value = 1107
$ws.Cells.Item(226, 20).Value = 7

# Make sure the date cells keep the same date/time display format as the
# rest of column D.
$ws.Range("D224:D226").NumberFormat = $ws.Range("D227").NumberFormat

$ws.Range("A1").Select()
